$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark. In the original document it
#    sits right at the end of the first paragraph (after the long
#    descriptive text about z-coordinate accuracy).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the "Cluster73" heading paragraph by its text (robust to
#    paragraph-index drift), then insert a brand-new, empty Heading3
#    paragraph immediately before it -- i.e. right after the paragraph
#    that contains the second figure (the "Cluster71" image).
$findRange = $d.Content
$null = $findRange.Find.Execute("Cluster73", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($findRange.Start -ge $p.Range.Start) -and ($findRange.Start -lt $p.Range.End)) {
        $targetIndex = $i
        break
    }
}

$cluster73 = $d.Paragraphs.Item($targetIndex)
$cluster73.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Style = "Heading 3"

# 3. Re-create the "_GoBack" bookmark collapsed inside that new, empty
#    paragraph (mirroring where Word leaves the bookmark after the
#    user's last edit location). A zero-length range placed directly at
#    an empty paragraph's start is mishandled by this host, so we work
#    around it: insert a placeholder character, wrap the bookmark around
#    it, then delete the placeholder -- the bookmark collapses cleanly
#    to zero length in place, exactly as Word itself would leave it.
$newPara.Range.InsertBefore("X")
$placeholder = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholder)
$placeholder2 = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$placeholder2.Delete()
